# Add a faculty-sheet header row: FSN | Name | Age | Phone
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "FSN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"

# Move the active selection to E1, matching where Excel leaves the cursor
# after typing the last header (Tab/Enter moves one cell to the right).
$ws.Range("E1").Select() | Out-Null
